$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 13): actuator test case with parabolic inflow ---
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "actuator_unsteady05"
$ws.Range("D13").Value = "u=parabolic"
$ws.Range("E13").Value = 100
$ws.Range("F13").Value = "80x40"
$ws.Range("G13").Value = "inflow-outflow; symmetry"
# Force "0.02" to be stored as text (matches dt column convention elsewhere in the table)
$ws.Range("H13").Value = "'0.02"
$ws.Range("H13").Style = "Normal"
$ws.Range("I13").Value = 10
$ws.Range("J13").Value = "RK44P2"
$ws.Range("K13").Value = "FOM"

# Highlight the copied boilerplate settings (Volumes..FOM/ROM) with a themed
# top border + light fill, same as Excel applies when a new table row is
# typed in under the existing banded style.
$band = $ws.Range("F13:K13")
$topBorder = $band.Borders.Item(8)
$topBorder.Color = 14461583
$topBorder.Weight = 2
$topBorder.LineStyle = 1
$band.Interior.Color = 15983578

# Trailing row outside the table, just like row 12 was before this edit
$ws.Range("A14").Value = 11

# --- Resize Table1 to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:O14"))

# --- Sheet view: selection + scroll position ---
$ws.Range("O13").Select()
